$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1529411764705882
$ws.Range("C2").Value = 0.6470588235294118
$ws.Range("J2").Value = 0.007843137254901961
$ws.Range("P2").Value = 0.1176470588235294
$ws.Range("S2").Value = 0.07450980392156863
$ws.Range("B3").Value = 0.01764705882352941
$ws.Range("C3").Value = 0.03529411764705882
$ws.Range("J3").Value = 0.06470588235294118
$ws.Range("P3").Value = 0.7235294117647059
$ws.Range("S3").Value = 0.1588235294117647
$ws.Range("J4").Value = 0.08695652173913043
$ws.Range("P4").Value = 0.6739130434782609
$ws.Range("S4").Value = 0.2391304347826087
$ws.Range("B6").Value = 0.07834101382488479
$ws.Range("D6").Value = 0.009216589861751152
$ws.Range("F6").Value = 0.1105990783410138
$ws.Range("J6").Value = 0.2672811059907834
$ws.Range("O6").Value = 0.004608294930875576
$ws.Range("Q6").Value = 0.119815668202765
$ws.Range("R6").Value = 0.05990783410138249
$ws.Range("S6").Value = 0.3502304147465438
$ws.Range("B7").Value = 0.1208791208791209
$ws.Range("D7").Value = 0.01648351648351648
$ws.Range("E7").Value = 0.005494505494505495
$ws.Range("F7").Value = 0.05494505494505494
$ws.Range("J7").Value = 0.08791208791208792
$ws.Range("O7").Value = 0.005494505494505495
$ws.Range("Q7").Value = 0.1648351648351648
$ws.Range("R7").Value = 0.1043956043956044
$ws.Range("S7").Value = 0.4395604395604396
$ws.Range("B8").Value = 0.0960591133004926
$ws.Range("D8").Value = 0.02955665024630542
$ws.Range("F8").Value = 0.08866995073891626
$ws.Range("J8").Value = 0.1280788177339902
$ws.Range("O8").Value = 0.01231527093596059
$ws.Range("Q8").Value = 0.1650246305418719
$ws.Range("R8").Value = 0.09359605911330049
$ws.Range("S8").Value = 0.3866995073891626
$ws.Range("B9").Value = 0.09195402298850575
$ws.Range("D9").Value = 0.02298850574712644
$ws.Range("F9").Value = 0.07471264367816093
$ws.Range("J9").Value = 0.1206896551724138
$ws.Range("O9").Value = 0.02298850574712644
$ws.Range("Q9").Value = 0.1494252873563219
$ws.Range("R9").Value = 0.1149425287356322
$ws.Range("S9").Value = 0.4022988505747127
$ws.Range("B10").Value = 0.09920983318700614
$ws.Range("D10").Value = 0.02458296751536436
$ws.Range("F10").Value = 0.05706760316066725
$ws.Range("J10").Value = 0.1255487269534679
$ws.Range("O10").Value = 0.01755926251097454
$ws.Range("Q10").Value = 0.2019315188762072
$ws.Range("R10").Value = 0.1115013169446883
$ws.Range("S10").Value = 0.3625987708516242
$ws.Range("G11").Value = 0.1388888888888889
$ws.Range("J11").Value = 0.06349206349206349
$ws.Range("K11").Value = 0.1706349206349206
$ws.Range("L11").Value = 0.6111111111111112
$ws.Range("S11").Value = 0.01587301587301587
$ws.Range("G12").Value = 0.7530864197530864
$ws.Range("J12").Value = 0.1975308641975309
$ws.Range("K12").Value = 0.006172839506172839
$ws.Range("L12").Value = 0.03703703703703703
$ws.Range("S12").Value = 0.006172839506172839
$ws.Range("F13").Value = 0.02439024390243903
$ws.Range("G13").Value = 0.6829268292682927
$ws.Range("J13").Value = 0.2439024390243902
$ws.Range("S13").Value = 0.04878048780487805
$ws.Range("F15").Value = 0.02727272727272727
$ws.Range("H15").Value = 0.1363636363636364
$ws.Range("I15").Value = 0.06363636363636363
$ws.Range("J15").Value = 0.3681818181818182
$ws.Range("K15").Value = 0.06818181818181818
$ws.Range("M15").Value = 0.00909090909090909
$ws.Range("O15").Value = 0.1045454545454545
$ws.Range("S15").Value = 0.2227272727272727
$ws.Range("F16").Value = 0.01136363636363636
$ws.Range("H16").Value = 0.1079545454545455
$ws.Range("I16").Value = 0.1022727272727273
$ws.Range("J16").Value = 0.4375
$ws.Range("K16").Value = 0.09659090909090909
$ws.Range("M16").Value = 0.03409090909090909
$ws.Range("N16").Value = 0.01136363636363636
$ws.Range("O16").Value = 0.09090909090909091
$ws.Range("S16").Value = 0.1079545454545455
$ws.Range("F17").Value = 0.02110817941952507
$ws.Range("H17").Value = 0.1846965699208443
$ws.Range("I17").Value = 0.09498680738786279
$ws.Range("J17").Value = 0.3931398416886543
$ws.Range("K17").Value = 0.1134564643799472
$ws.Range("M17").Value = 0.01319261213720317
$ws.Range("O17").Value = 0.0870712401055409
$ws.Range("S17").Value = 0.09234828496042216
$ws.Range("F18").Value = 0.02803738317757009
$ws.Range("H18").Value = 0.1775700934579439
$ws.Range("I18").Value = 0.08878504672897196
$ws.Range("J18").Value = 0.4158878504672897
$ws.Range("K18").Value = 0.1168224299065421
$ws.Range("M18").Value = 0.009345794392523364
$ws.Range("N18").Value = 0.004672897196261682
$ws.Range("O18").Value = 0.05607476635514019
$ws.Range("S18").Value = 0.102803738317757
$ws.Range("F19").Value = 0.01914311759343664
$ws.Range("H19").Value = 0.2297174111212398
$ws.Range("I19").Value = 0.0829535095715588
$ws.Range("J19").Value = 0.3555150410209663
$ws.Range("K19").Value = 0.097538742023701
$ws.Range("M19").Value = 0.02734731084776664
$ws.Range("N19").Value = 0.004557885141294439
$ws.Range("O19").Value = 0.07292616226071102
$ws.Range("S19").Value = 0.1103008204193254
